$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.98
$ws.Range("C2").Value = 9.4600000000000009

$ws.Range("B3").Value = 9.86
$ws.Range("C3").Value = 75.069999999999993

$ws.Range("B4").Value = 10.01
$ws.Range("C4").Value = 69.09
